$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.169409036636353
$ws.Range("B1").Value = 2.926520586013794
$ws.Range("C1").Value = 4.455004215240479
$ws.Range("D1").Value = 1.958391785621643
$ws.Range("E1").Value = 1.161061525344849
